# NSMB - Starting 8-2
# Append the tracking rows for the 8-2 segment under the "World 8" table on
# the V4 sheet (rows 92-98), continuing directly after the existing
# "Enter 8-1" / "1st Move" rows (90-91).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# Label, Known(B), TAS(C) for each new row, in sheet order.
$rows = @(
    @{ Row = 92; Label = "Checkpoint 2005";    Known = 27730; TAS = 32556 },
    @{ Row = 93; Label = "Get flag";            Known = 28335; TAS = 33163 },
    @{ Row = 94; Label = "End level";           Known = 28853; TAS = 33681 },
    @{ Row = 95; Label = "Enter 8-2";           Known = 29213; TAS = 34468 },
    @{ Row = 96; Label = "1st Move";            Known = 29438; TAS = 34713 },
    @{ Row = 97; Label = "Enter pipe";          Known = 29585; TAS = 34860 },
    @{ Row = 98; Label = "Mario touch ground";  Known = 29708; TAS = 34998 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value2 = $r.Label
    $ws.Cells.Item($rowNum, 2).Value2 = $r.Known
    $ws.Cells.Item($rowNum, 3).Value2 = $r.TAS
    $ws.Cells.Item($rowNum, 4).Formula = "=IF(B" + $rowNum + ">0,C" + $rowNum + "-B" + $rowNum + ",0)"

    # Match the formatting already used by the rest of the table (rows 90-91).
    $ws.Cells.Item($rowNum, 1).Style = $ws.Cells.Item(90, 1).Style
    $ws.Cells.Item($rowNum, 2).Style = $ws.Cells.Item(90, 2).Style
    $ws.Cells.Item($rowNum, 3).Style = $ws.Cells.Item(90, 3).Style
    $ws.Cells.Item($rowNum, 4).Style = $ws.Cells.Item(90, 4).Style
}

# Keep the frozen pane / selection in sync with the new bottom of the table,
# same as Excel does after scrolling to and selecting the next empty row.
$ws.Application.Goto($ws.Range("A99"))
$ws.Range("A99").Select()
